# "Update countries & provincias Spain"
#
# Refreshes the Spain provinces COVID table with a newer data snapshot
# (21:28 instead of 18:44). The sheet is kept sorted descending by
# "Casos totales" (column B), so several rows swap places as a result of
# the new counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "last refreshed" timestamp banner in A1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 27 de Marzo de 2020 a las 21:28"

function Set-Row($row, $name, $total, $active, $recovered, $deaths) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $total
    $ws.Cells.Item($row, 3).Value = $active
    $ws.Cells.Item($row, 4).Value = $recovered
    $ws.Cells.Item($row, 5).Value = $deaths
}

# Rows 18-20: Aragon/Asturias/Gran Canaria re-sort with refreshed figures
Set-Row 18 "Asturias" 1004 65 906 33
Set-Row 19 "Gran Canaria" 914 20 866 28
Set-Row 20 "Aragon" 907 29 838 40

# Rows 24-28: Murcia/Pontevedra/Albacete/Granada/Sevilla re-sort with refreshed figures
Set-Row 24 "Murcia" 800 12 771 17
Set-Row 25 "Pontevedra" 795 67 769 8
Set-Row 26 "Albacete" 780 153 667 83
Set-Row 27 "Granada" 711 1 676 34
Set-Row 28 "Sevilla" 708 8 675 25

# Deaths ("Muertes") count updated from 27 to 28 for several other provinces
$ws.Cells.Item(31, 5).Value = 28
$ws.Cells.Item(55, 5).Value = 28
$ws.Cells.Item(57, 5).Value = 28
$ws.Cells.Item(58, 5).Value = 28
$ws.Cells.Item(62, 5).Value = 28
$ws.Cells.Item(63, 5).Value = 28
